$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 213
$ws.Range("I11").Value = 213
$ws.Range("K11").Value = 213
$ws.Range("M11").Value = -73
$ws.Range("H12").Value = 926.8570999999999
$ws.Range("I12").Value = 1049.3334
$ws.Range("J12").Value = 192
$ws.Range("K12").Value = 1049.3334
$ws.Range("L12").Value = 192
$ws.Range("M12").Value = -879.3334
$ws.Range("N12").Value = -532
$ws.Range("H28").Value = 1091.2
$ws.Range("I28").Value = 922
$ws.Range("J28").Value = 1260.4
$ws.Range("K28").Value = 922
$ws.Range("L28").Value = 1260.4
$ws.Range("M28").Value = -437
$ws.Range("N28").Value = -2230.4
$ws.Range("H41").Value = 415.9
$ws.Range("J41").Value = 1150
$ws.Range("L41").Value = 1150
$ws.Range("N41").Value = -2030
$ws.Range("H106").Value = 4347.5
$ws.Range("I106").Value = 4347.5
$ws.Range("K106").Value = 4347.5
$ws.Range("M106").Value = -3716.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3114.4
$ws.Range("I2").Value = 1393.5
$ws.Range("K2").Value = 1393.5
$ws.Range("M2").Value = -1280.5
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 18.5
$ws.Range("K5").Value = 5.5
$ws.Range("L5").Value = 18.5
$ws.Range("M5").Value = 106.5
$ws.Range("N5").Value = -242.5
$ws.Range("H31").Value = 9499.5
$ws.Range("I31").Value = 9499.5
$ws.Range("K31").Value = 9499.5
$ws.Range("M31").Value = -9205.5
$ws.Range("H61").Value = 936.5
$ws.Range("I61").Value = 936.5
$ws.Range("K61").Value = 936.5
$ws.Range("M61").Value = -724.5
$ws.Range("H110").Value = 851
$ws.Range("I110").Value = 851
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 851
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1194
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 3114.4
$ws.Range("I116").Value = 1393.5
$ws.Range("K116").Value = 1393.5
$ws.Range("M116").Value = 900.5
$ws.Range("H122").Value = 1863.2858
$ws.Range("I122").Value = 1840.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5521.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3071.5
$ws.Range("N122").Value = -10900
$ws.Range("H136").Value = 936.5
$ws.Range("I136").Value = 936.5
$ws.Range("K136").Value = 2809.5
$ws.Range("M136").Value = -259.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3114.4
$ws.Range("I3").Value = 1393.5
$ws.Range("K3").Value = 1393.5
$ws.Range("M3").Value = -1279.5
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 5.5
$ws.Range("J4").Value = 18.5
$ws.Range("K4").Value = 5.5
$ws.Range("L4").Value = 18.5
$ws.Range("M4").Value = 109.5
$ws.Range("N4").Value = -248.5
$ws.Range("H8").Value = 80
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H64").Value = 1001.25
$ws.Range("J64").Value = 1003.5
$ws.Range("L64").Value = 1003.5
$ws.Range("N64").Value = -1453.5
$ws.Range("H67").Value = 1001.25
$ws.Range("J67").Value = 1003.5
$ws.Range("L67").Value = 1003.5
$ws.Range("N67").Value = -2563.5
$ws.Range("H86").Value = 1680.4445
$ws.Range("I86").Value = 2229
$ws.Range("J86").Value = 583.3333
$ws.Range("K86").Value = 2229
$ws.Range("L86").Value = 583.3333
$ws.Range("M86").Value = -1106
$ws.Range("N86").Value = -2829.3333
$ws.Range("H89").Value = 1680.4445
$ws.Range("I89").Value = 2229
$ws.Range("J89").Value = 583.3333
$ws.Range("K89").Value = 11145
$ws.Range("L89").Value = 2916.6665
$ws.Range("M89").Value = -5529
$ws.Range("N89").Value = -14148.6665
$ws.Range("H94").Value = 1719.9231
$ws.Range("I94").Value = 1896.2727
$ws.Range("K94").Value = 1896.2727
$ws.Range("M94").Value = -1445.2727
$ws.Range("H107").Value = 669.85187
$ws.Range("I107").Value = 668.6923
$ws.Range("K107").Value = 668.6923
$ws.Range("M107").Value = 1251.3077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4999.3335
$ws.Range("H16").Value = 2400
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 6000
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -913
$ws.Range("N16").Value = -6574
$ws.Range("H58").Value = 3235.5
$ws.Range("I58").Value = 3235.5
$ws.Range("K58").Value = 3235.5
$ws.Range("M58").Value = -3032.5
$ws.Range("H113").Value = 2400
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -10340
$ws.Range("H136").Value = 3235.5
$ws.Range("I136").Value = 3235.5
$ws.Range("K136").Value = 9706.5
$ws.Range("M136").Value = -7156.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100339.7
$ws.Range("I4").Value = 377.44446
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 1132.33338
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -1020.33338
$ws.Range("N4").Value = -3000224
$ws.Range("H12").Value = 97.625
$ws.Range("J12").Value = 97.625
$ws.Range("L12").Value = 292.875
$ws.Range("N12").Value = -638.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4399.6
$ws.Range("I70").Value = 4399.6
$ws.Range("K70").Value = 4399.6
$ws.Range("M70").Value = -4129.6
$ws.Range("H73").Value = 4399.6
$ws.Range("I73").Value = 4399.6
$ws.Range("K73").Value = 4399.6
$ws.Range("M73").Value = -3463.6
$ws.Range("H102").Value = 2140
$ws.Range("I102").Value = 1913.3334
$ws.Range("K102").Value = 1913.3334
$ws.Range("M102").Value = -291.3334
$ws.Range("H113").Value = 5101.75
$ws.Range("I113").Value = 4190.5557
$ws.Range("J113").Value = 7835.3335
$ws.Range("K113").Value = 4190.5557
$ws.Range("L113").Value = 7835.3335
$ws.Range("M113").Value = -2020.5557
$ws.Range("N113").Value = -12175.3335
$ws.Range("H132").Value = 2270.4614
$ws.Range("J132").Value = 2334
$ws.Range("L132").Value = 7002
$ws.Range("N132").Value = -12062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 400
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H16").Value = 264.83334
$ws.Range("I16").Value = 257.8
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 257.8
$ws.Range("L16").Value = 300
$ws.Range("M16").Value = -87.80000000000001
$ws.Range("N16").Value = -640
$ws.Range("H22").Value = 2016
$ws.Range("I22").Value = 456.14285
$ws.Range("J22").Value = 4199.8
$ws.Range("K22").Value = 456.14285
$ws.Range("L22").Value = 4199.8
$ws.Range("M22").Value = -161.14285
$ws.Range("N22").Value = -4789.8
$ws.Range("H27").Value = 2016
$ws.Range("I27").Value = 456.14285
$ws.Range("J27").Value = 4199.8
$ws.Range("K27").Value = 456.14285
$ws.Range("L27").Value = 4199.8
$ws.Range("M27").Value = -349.14285
$ws.Range("N27").Value = -4413.8
$ws.Range("H46").Value = 4453.885
$ws.Range("H61").Value = 3999.5
$ws.Range("J61").Value = 3999.5
$ws.Range("L61").Value = 3999.5
$ws.Range("N61").Value = -4403.5
$ws.Range("H82").Value = 700
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 800
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = -1522
$ws.Range("H85").Value = 700
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 800
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 800
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = -3296
$ws.Range("H100").Value = 22100
$ws.Range("I100").Value = 10000
$ws.Range("K100").Value = 10000
$ws.Range("M100").Value = -9459
$ws.Range("H113").Value = 3999.5
$ws.Range("J113").Value = 3999.5
$ws.Range("L113").Value = 3999.5
$ws.Range("N113").Value = -8339.5
$ws.Range("H122").Value = 4373.25
$ws.Range("I122").Value = 4373.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13119.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10669.75
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 667707.7
$ws.Range("I3").Value = 667707.7
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 667707.7
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -667593.7
$ws.Range("N3").ClearContents()
$ws.Range("H96").Value = 1735.1177
$ws.Range("J96").Value = 1619.5
$ws.Range("L96").Value = 1619.5
$ws.Range("N96").Value = -4365.5
$ws.Range("H113").Value = 526.94116
$ws.Range("I113").Value = 647.7143
$ws.Range("J113").Value = 442.4
$ws.Range("K113").Value = 1943.1429
$ws.Range("L113").Value = 1327.2
$ws.Range("M113").Value = 226.8571000000002
$ws.Range("N113").Value = -5667.2
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H126").Value = 1833.3334
$ws.Range("I126").Value = 1833.3334
$ws.Range("K126").Value = 5500.0002
$ws.Range("M126").Value = -3030.0002
